# Fruta / hortaliza, semanal
# Insert a new weekly observation row at row 335 (pushing the existing
# rows 335-350 down to 336-351), then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 335:350 down to 336:351, duplicating row 335's formatting
# (style s="2" on column D) into the newly inserted row - same as Excel's
# native "Insert Copied/Sheet Rows" behaviour.
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new weekly record.
$ws.Range("A335").Value = 4
$ws.Range("B335").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C335").Value = "Los Lagos"
$ws.Range("D335").Value = 45041
$ws.Range("E335").Value = 10
$ws.Range("F335").Value = 100112039
$ws.Range("G335").Value = "Ciboulette"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 240
$ws.Range("K335").Value = 3500
$ws.Range("L335").Value = 3500
$ws.Range("M335").Value = 3500
$ws.Range("N335").Value = "`$/docena de atados"
$ws.Range("O335").Value = "Región Metropolitana"
$ws.Range("P335").Value = 1167
$ws.Range("Q335").Value = 3
$ws.Range("R335").Value = "Hortaliza"
